# Auto detect the starting (first empty) row on the active ("total_data") sheet
# and append a new record there.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A, then the next row is where new data starts.
$lastRow = $ws.Cells($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Make sure the new cells use the same "text" number format as the rest of the
# data rows so values are stored as shared strings (matching existing rows).
$rng = $ws.Range($ws.Cells($newRow, 1), $ws.Cells($newRow, 4))
$rng.NumberFormat = "@"

$ws.Cells($newRow, 1).Value = "31"
$ws.Cells($newRow, 2).Value = "May"
$ws.Cells($newRow, 3).Value = "Bach"
$ws.Cells($newRow, 4).Value = "Laos"
